# Update the "具体时间范围" text in E2/E3 on the sheets that contain this data.
# The diff removes the spaces around the dash: "09:00 - 05." -> "09:00-05."

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("E2").Value = "2024.05.01 09:00-05.01 17:00"
    $ws.Range("E3").Value = "2024.05.18 09:00-05.18 17:00"
}
